$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# NOTE: the order in which brand-new string values are first assigned
# controls the order they are appended to the shared-string table, so we
# seed them in the exact order required: "fail", "Fail", "WIN", "15290 possible".
$ws.Cells.Item(35, 11).Value = "fail"
$ws.Cells.Item(42, 11).Value = "Fail"
$ws.Cells.Item(52, 11).Value = "WIN"
$ws.Cells.Item(39, 9).Value = "15290 possible"

# --- Row 33 / 34: only the row span widens (handled implicitly by adding
#     cells further down the sheet); nothing else to change here.

# --- New helper columns I/J/K/L (rows 35-52) + updated B/H values ---

# Row 35
$ws.Cells.Item(35, 10).Value = 76

# Row 36
$ws.Cells.Item(36, 10).Value = 77
$ws.Cells.Item(36, 11).Value = "fail"

# Row 37
$ws.Cells.Item(37, 10).Value = 78

# Row 38
$ws.Cells.Item(38, 10).Value = 79
$ws.Cells.Item(38, 11).Value = "fail"

# Row 39 - also update B39, add J39
$ws.Range("B39").Value = 15307
$ws.Cells.Item(39, 10).Value = 80

# Row 40 - update B40, add J40
$ws.Range("B40").Value = 15573
$ws.Cells.Item(40, 10).Value = 81

# Row 41
$ws.Cells.Item(41, 10).Value = 82

# Row 42 - update B42, add H42, J42
$ws.Range("B42").Value = 15722
$ws.Range("H42").Value = 15781
$ws.Cells.Item(42, 10).Value = 83

# Row 43
$ws.Cells.Item(43, 10).Value = 84

# Row 44
$ws.Cells.Item(44, 10).Value = 85
$ws.Cells.Item(44, 11).Value = "Fail"

# Row 45
$ws.Cells.Item(45, 10).Value = 86

# Row 46
$ws.Cells.Item(46, 10).Value = 87

# Row 47
$ws.Cells.Item(47, 10).Value = 88

# Row 48
$ws.Cells.Item(48, 10).Value = 89

# Row 49
$ws.Cells.Item(49, 10).Value = 90

# Row 50
$ws.Cells.Item(50, 10).Value = 91

# Row 51
$ws.Cells.Item(51, 10).Value = 92
$ws.Cells.Item(51, 11).Value = "fail"

# Row 52
$ws.Cells.Item(52, 9).Value = 95090
$ws.Cells.Item(52, 10).Value = 93
$ws.Cells.Item(52, 12).Value = 94790

# --- sheetView pane/selection changes ---
$window = $ws.Application.ActiveWindow
$window.ScrollRow = 18
$ws.Range("K25").Select()
